$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New OHLC / shares_outstanding / fixed_ticker values for each dated row (2-30).
# The source data previously mixed in rows for other tickers (FTNT, INTC, TRI, ...);
# this consolidates every row onto the DDOG ticker with its correct price history,
# removing the need for those extra shared strings.
$rowData = @(
    @{ Row = 2; D = 40.34999847412109; E = 33.90999984741211; F = 41.43999862670898; G = 30.46999931335449; H = 323270704 }
    @{ Row = 3; D = 40.34999847412109; E = 33.90999984741211; F = 41.43999862670898; G = 30.46999931335449; H = 323270704 }
    @{ Row = 4; D = 40.34999847412109; E = 33.90999984741211; F = 41.43999862670898; G = 30.46999931335449; H = 323270704 }
    @{ Row = 5; D = 40.34999847412109; E = 33.90999984741211; F = 41.43999862670898; G = 30.46999931335449; H = 323270704 }
    @{ Row = 6; D = 40.34999847412109; E = 33.90999984741211; F = 41.43999862670898; G = 30.46999931335449; H = 323270704 }
    @{ Row = 7; D = 34.25; E = 33.59000015258789; F = 37.97000122070312; G = 27.54999923706055; H = 323270704 }
    @{ Row = 8; D = 38.22000122070312; E = 46.20999908447266; F = 47.2599983215332; G = 35.95800018310547; H = 323270704 }
    @{ Row = 9; D = 34.83000183105469; E = 45.11999893188477; F = 45.27999877929688; G = 33.04000091552734; H = 323270704 }
    @{ Row = 10; D = 87.19000244140625; E = 93.86000061035156; F = 98.98999786376952; G = 81.54100036621094; H = 323270704 }
    @{ Row = 11; D = 101.4300003051758; E = 90.75; F = 118.129997253418; G = 90.54000091552734; H = 323270704 }
    @{ Row = 12; D = 98.69000244140624; E = 102.75; F = 112.2900009155273; G = 89.83200073242188; H = 323270704 }
    @{ Row = 13; D = 86.94999694824219; E = 85.7699966430664; F = 96.3499984741211; G = 81.69999694824219; H = 323270704 }
    @{ Row = 14; D = 103.5100021362305; E = 110.6999969482422; F = 114.1999969482422; G = 101.379997253418; H = 323270704 }
    @{ Row = 15; D = 142.0200042724609; E = 167.0500030517578; F = 168.7299957275391; G = 134.6300048828125; H = 323270704 }
    @{ Row = 16; D = 178.75; E = 146.1100006103516; F = 180.2799987792969; G = 119.1900024414062; H = 323270704 }
    @{ Row = 17; D = 150.9799957275391; E = 120.7799987792969; F = 159; G = 117.8600006103516; H = 323270704 }
    @{ Row = 18; D = 95.25; E = 102.0100021362305; F = 112.5800018310547; G = 84.45999908447266; H = 323270704 }
    @{ Row = 19; D = 89.61000061035156; E = 80.51000213623047; F = 97.31400299072266; G = 75.53500366210938; H = 323270704 }
    @{ Row = 20; D = 75.19499969482422; E = 74.80999755859375; F = 78.62999725341797; G = 61.34000015258789; H = 323270704 }
    @{ Row = 21; D = 71.61499786376953; E = 67.37999725341797; F = 72.66999816894531; G = 62.59700012207031; H = 323270704 }
    @{ Row = 22; D = 98.0999984741211; E = 116.7200012207031; F = 118.0199966430664; G = 95.30999755859376; H = 323270704 }
    @{ Row = 23; D = 92.54000091552734; E = 81.47000122070312; F = 94.8000030517578; G = 79.37000274658203; H = 323270704 }
    @{ Row = 24; D = 119.2249984741211; E = 124.4400024414062; F = 136.0749969482422; G = 110.6949996948242; H = 323270704 }
    @{ Row = 25; D = 123.8649978637695; E = 125.5; F = 132.8999938964844; G = 118.4599990844727; H = 323270704 }
    @{ Row = 26; D = 129.0399932861328; E = 116.4400024414062; F = 134.8000030517578; G = 113.5800018310547; H = 323270704 }
    @{ Row = 27; D = 117.4499969482422; E = 125.4400024414062; F = 131.8399963378906; G = 113.0800018310547; H = 323270704 }
    @{ Row = 28; D = 145; E = 142.7100067138672; F = 153.4199981689453; G = 134.7149963378906; H = 323270704 }
    @{ Row = 29; D = 100.1900024414062; E = 102.1600036621094; F = 102.8600006103516; G = 81.62999725341797; H = 323270704 }
    @{ Row = 30; D = 133.6600036621094; E = 139.9799957275391; F = 157.375; G = 130.2799987792969; H = 323270704 }
)

foreach ($item in $rowData) {
    $r = $item.Row
    $ws.Range("D$r").Value = $item.D
    $ws.Range("E$r").Value = $item.E
    $ws.Range("F$r").Value = $item.F
    $ws.Range("G$r").Value = $item.G
    $ws.Range("H$r").Value = $item.H
    $ws.Range("I$r").Value = "DDOG"
}